# "State Mapping Status" workbook update
# Updates a handful of state-contact-tracking rows on the "SEARCH" sheet:
#   - Wisconsin (row 56): new "Mapping Received" note plus a note about the
#     located DOJ statute crosswalk table
#   - Oregon (row 52): new note about contacting Kelly Officer / CJIS
#   - New Mexico (row 51): "Existing Mapping" status updated, new "Mapping
#     Received" note added, and the notes column rewritten with the latest
#     status from NM DPS / the Statistical Analysis Center
#   - Mississippi (row 50): contact name updated, new note added

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wisconsin
$ws.Range("F56").Value = "Found via research"
$ws.Range("J56").Value = "Located Wisconsin DOJ statute table crosswalk to NCIC.pdf.  Includes Statute, Offense Code, Class, title, NCIC.  Cathy will upload to Github."

# Oregon
$ws.Range("J52").Value = "Spoke with Kelly Officer, and CJIS folks.  Submitted a written request to Michael Hawkins (CCH/LEDS Program Manager) for OR statute ot UCR coding."

# New Mexico
$ws.Range("E51").Value = "No"
$ws.Range("F51").Value = "Found via research"
$ws.Range("J51").Value = "Tim responded that NM DPS did not have a table. Cathy followed up with the Statistical Analysis Center; a consortium of agncies, including NM DPS, has created a common charge code table"

# Mississippi
$ws.Range("J50").Value = "Wrote to Captain Wilson requesting mapping artifacts"
$ws.Range("I50").Value = "Lt. Charlie Hill/Capt. Wilson"
